$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 63

# Duplicate the formatting of the previous (last existing) data row so that the
# new row gets the same cell styles: bold/bordered index column (A) and the
# date/time number format on the match-date column (E).
$ws.Range("A62:V62").Copy()
$ws.Range("A63:V63").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 62
$ws.Cells.Item($row, 2).Value = "croatia"
$ws.Cells.Item($row, 3).Value = "hnl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45233.75
$ws.Cells.Item($row, 6).Value = "Lok. Zagreb"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Rijeka"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 3.24
$ws.Cells.Item($row, 11).Value = "29/10/2023 17:13"
$ws.Cells.Item($row, 12).Value = 4.19
$ws.Cells.Item($row, 13).Value = "03/11/2023 17:55"
$ws.Cells.Item($row, 14).Value = 3.19
$ws.Cells.Item($row, 15).Value = "29/10/2023 17:13"
$ws.Cells.Item($row, 16).Value = 3.45
$ws.Cells.Item($row, 17).Value = "03/11/2023 17:55"
$ws.Cells.Item($row, 18).Value = 2.23
$ws.Cells.Item($row, 19).Value = "29/10/2023 17:13"
$ws.Cells.Item($row, 20).Value = 1.93
$ws.Cells.Item($row, 21).Value = "03/11/2023 17:51"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/lok-zagreb-rijeka/rydIYL4m/"
